# Insert two new data rows right after the existing row 1019 (i.e. at row 1020),
# shifting every row from the old 1020 onward down by two. Then populate the
# two freshly-inserted rows with the new "Terminal Hortofrutícola Agro Chillán"
# apple-price records (Fuji royal, Primera + Segunda, fecha 44931, Región de
# O'Higgins).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 1020.
$ws.Rows.Item(1020).EntireRow.Insert()
$ws.Rows.Item(1020).EntireRow.Insert()

# Constant columns shared by every record on this sheet.
$mercadoId = 7
$mercado   = "Terminal Hortofrutícola Agro Chillán"
$region    = "Ñuble"
$codreg    = 16
$tipo      = "Fruta"
$productoId = 100104
$producto  = "Frutos de pepita"
$categoriaId = 100104002
$categoria = "Manzana"
$unidad    = "`$/caja 16 kilos empedrada"
$kgUnidad  = 16

function Set-Record {
    param($Row, $Fecha, $Variedad, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $Origen, $PrecioKg)

    $ws.Cells.Item($Row, 1).Value  = $mercadoId
    $ws.Cells.Item($Row, 2).Value  = $mercado
    $ws.Cells.Item($Row, 3).Value  = $region
    $ws.Cells.Item($Row, 4).Value  = $Fecha
    $ws.Cells.Item($Row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($Row, 5).Value  = $codreg
    $ws.Cells.Item($Row, 6).Value  = $tipo
    $ws.Cells.Item($Row, 7).Value  = $productoId
    $ws.Cells.Item($Row, 8).Value  = $producto
    $ws.Cells.Item($Row, 9).Value  = $categoriaId
    $ws.Cells.Item($Row, 10).Value = $categoria
    $ws.Cells.Item($Row, 11).Value = $Variedad
    $ws.Cells.Item($Row, 12).Value = $Calidad
    $ws.Cells.Item($Row, 13).Value = $Volumen
    $ws.Cells.Item($Row, 14).Value = $PrecioMin
    $ws.Cells.Item($Row, 15).Value = $PrecioMax
    $ws.Cells.Item($Row, 16).Value = $PrecioProm
    $ws.Cells.Item($Row, 17).Value = $unidad
    $ws.Cells.Item($Row, 18).Value = $Origen
    $ws.Cells.Item($Row, 19).Value = $PrecioKg
    $ws.Cells.Item($Row, 20).Value = $kgUnidad
}

# New row 1020: Fuji royal, Primera, fecha 2023-01-05 (44931), Región de O'Higgins
Set-Record 1020 44931 "Fuji royal" "Primera" 120 19000 20000 19500 "Región de O'Higgins" 1219

# New row 1021: Fuji royal, Segunda, fecha 2023-01-05 (44931), Región de O'Higgins
Set-Record 1021 44931 "Fuji royal" "Segunda" 80 18000 18000 18000 "Región de O'Higgins" 1125
